$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell values that look like plain decimal numbers must be forced to text
# (matching the source data which stores prices/percentages as inline strings)
# by temporarily applying a text number format, then resetting the style so
# no extra formatting is left behind on the cell.

$ws.Range('D2').Value = '42.863.33'
$ws.Range('E2').Value = '  -0.29%  '
$ws.Range('D3').Value = '2.298.40'
$ws.Range('E3').Value = '  -0.02%  '
$ws.Range('E4').Value = '  -0.04%  '
$c = $ws.Range('D5')
$c.NumberFormat = "@"
$c.Value = '305.17'
$c.Style = "Normal"
$ws.Range('E5').Value = '  +1.56%  '
$c = $ws.Range('D6')
$c.NumberFormat = "@"
$c.Value = '97.39'
$c.Style = "Normal"
$ws.Range('E6').Value = '  +0.07%  '
$ws.Range('E7').Value = '  -1.58%  '
$ws.Range('E8').Value = '  -0.03%  '
$c = $ws.Range('D9')
$c.NumberFormat = "@"
$c.Value = '0.505'
$c.Style = "Normal"
$ws.Range('E9').Value = '  -1.75%  '
$c = $ws.Range('D10')
$c.NumberFormat = "@"
$c.Value = '35.62'
$c.Style = "Normal"
$ws.Range('E10').Value = '  +0.26%  '
$ws.Range('E11').Value = '  +0.07%  '
$c = $ws.Range('D12')
$c.NumberFormat = "@"
$c.Value = '18.18'
$c.Style = "Normal"
$ws.Range('E12').Value = '  +1.46%  '
$ws.Range('E13').Value = '  +1.17%  '
$c = $ws.Range('D14')
$c.NumberFormat = "@"
$c.Value = '6.78'
$c.Style = "Normal"
$ws.Range('E14').Value = '  -0.92%  '
$ws.Range('D15').Value = '2.657.77'
$ws.Range('E15').Value = '  -0.12%  '
$ws.Range('D16').Value = '2.300.59'
$ws.Range('E16').Value = '  -0.35%  '
$ws.Range('E17').Value = '  +0.03%  '
$ws.Range('D18').Value = '42.812.86'
$ws.Range('E18').Value = '  -0.23%  '
$c = $ws.Range('D19')
$c.NumberFormat = "@"
$c.Value = '12.68'
$c.Style = "Normal"
$ws.Range('E19').Value = '  -4.96%  '
$ws.Range('E20').Value = '  -0.18%  '
$ws.Range('E21').Value = '  -1.04%  '
$c = $ws.Range('D22')
$c.NumberFormat = "@"
$c.Value = '67.87'
$c.Style = "Normal"
$ws.Range('E22').Value = '  -0.52%  '
$c = $ws.Range('D23')
$c.NumberFormat = "@"
$c.Value = '236.20'
$c.Style = "Normal"
$ws.Range('E23').Value = '  -0.82%  '
$c = $ws.Range('D24')
$c.NumberFormat = "@"
$c.Value = '2.14'
$c.Style = "Normal"
$ws.Range('E24').Value = '  -2.09%  '
$ws.Range('E25').Value = '  +1.80%  '
$ws.Range('E26').Value = '  +0.14%  '
$ws.Range('E27').Value = '  -0.13%  '
$c = $ws.Range('D28')
$c.NumberFormat = "@"
$c.Value = '25.35'
$c.Style = "Normal"
$ws.Range('E28').Value = '  +2.92%  '
$c = $ws.Range('D29')
$c.NumberFormat = "@"
$c.Value = '167.54'
$c.Style = "Normal"
$ws.Range('E29').Value = '  -0.60%  '
$c = $ws.Range('D30')
$c.NumberFormat = "@"
$c.Value = '2.05'
$c.Style = "Normal"
$ws.Range('E30').Value = '  +1.18%  '
$ws.Range('E31').Value = '  -0.85%  '
$c = $ws.Range('D32')
$c.NumberFormat = "@"
$c.Value = '33.16'
$c.Style = "Normal"
$ws.Range('E32').Value = '  +1.07%  '
$ws.Range('E33').Value = '  +0.02%  '
$c = $ws.Range('D34')
$c.NumberFormat = "@"
$c.Value = '4.80'
$c.Style = "Normal"
$ws.Range('E34').Value = '  +0.07%  '
$ws.Range('E35').Value = '  -2.69%  '
$c = $ws.Range('D36')
$c.NumberFormat = "@"
$c.Value = '17.25'
$c.Style = "Normal"
$ws.Range('E36').Value = '  -4.62%  '
$c = $ws.Range('D38')
$c.NumberFormat = "@"
$c.Value = '0.0691'
$c.Style = "Normal"
$ws.Range('E38').Value = '  +0.66%  '
$ws.Range('E39').Value = '  -1.31%  '
$ws.Range('E40').Value = '  -1.32%  '
$ws.Range('E41').Value = '  -1.14%  '
$ws.Range('E42').Value = '  -0.28%  '
$ws.Range('D43').Value = '2.005.41'
$ws.Range('E43').Value = '  -0.27%  '
$ws.Range('E44').Value = '  -2.14%  '
$c = $ws.Range('D45')
$c.NumberFormat = "@"
$c.Value = '10.02'
$c.Style = "Normal"
$ws.Range('E45').Value = '  -1.77%  '
$ws.Range('B46').Value = 'EnergySwap'
$ws.Range('C46').Value = 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens'
$c = $ws.Range('D46')
$c.NumberFormat = "@"
$c.Value = '17.83'
$c.Style = "Normal"
$ws.Range('E46').Value = '  +3.12%  '
$ws.Range('B47').Value = 'ApeXProtocol'
$ws.Range('C47').Value = 'https://coinranking.com/coin/ze0N2Rcyu+apexprotocol-apex'
$c = $ws.Range('D47')
$c.NumberFormat = "@"
$c.Value = '2.09'
$c.Style = "Normal"
$ws.Range('E47').Value = '  -3.35%  '
$c = $ws.Range('D48')
$c.NumberFormat = "@"
$c.Value = '2.78'
$c.Style = "Normal"
$ws.Range('E48').Value = '  -1.26%  '
$ws.Range('B49').Value = 'MultiversX'
$ws.Range('C49').Value = 'https://coinranking.com/coin/omwkOTglq+multiversx-egld'
$c = $ws.Range('D49')
$c.NumberFormat = "@"
$c.Value = '53.83'
$c.Style = "Normal"
$ws.Range('E49').Value = '  -0.50%  '
$ws.Range('B50').Value = 'HuobiToken'
$ws.Range('C50').Value = 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht'
$c = $ws.Range('D50')
$c.NumberFormat = "@"
$c.Value = '2.84'
$c.Style = "Normal"
$ws.Range('E50').Value = '  +1.02%  '
$ws.Range('D51').Value = '2.525.36'
$ws.Range('E51').Value = '  -0.05%  '
